# Update the imputed values for columns A and B (RandomForest algorithm
# result re-run) on Sheet1, as described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 8.325800000000003
$ws.Range("A3").Value = -22.02100000000001
$ws.Range("A14").Value = -21.77959999999999
$ws.Range("A16").Value = -21.64679999999998
$ws.Range("B18").Value = 6.342399999999994
$ws.Range("A21").Value = -19.93579999999998
$ws.Range("A23").Value = -20.96139999999998
$ws.Range("B24").Value = 6.399100000000001
$ws.Range("A25").Value = -21.49299999999998
$ws.Range("B25").Value = 5.9645
$ws.Range("A26").Value = -20.97849999999997
$ws.Range("B27").Value = 5.805100000000003
$ws.Range("A29").Value = -21.06559999999998
$ws.Range("B30").Value = 6.094899999999998
$ws.Range("B31").Value = 5.650400000000002
$ws.Range("B39").Value = 9.388500000000004
$ws.Range("A40").Value = -19.42389999999999
$ws.Range("B42").Value = 10.31199999999999
$ws.Range("B48").Value = 5.389800000000003
$ws.Range("B51").Value = 5.7369
$ws.Range("B52").Value = 5.575
$ws.Range("A53").Value = -21.8599
$ws.Range("B55").Value = 6.449199999999993
$ws.Range("B56").Value = 5.134499999999997
$ws.Range("A57").Value = -22.32530000000001
$ws.Range("B57").Value = 5.112799999999993
$ws.Range("A59").Value = -22.441
$ws.Range("B60").Value = 5.3826
$ws.Range("A65").Value = -21.84899999999998
$ws.Range("A69").Value = -21.6401
$ws.Range("B73").Value = 9.0976
$ws.Range("B74").Value = 9.780499999999993
$ws.Range("A79").Value = -20.67500000000002
$ws.Range("A83").Value = -21.76089999999999
$ws.Range("B89").Value = 4.589599999999994
$ws.Range("B90").Value = 5.621500000000001
$ws.Range("A91").Value = -21.44390000000002
$ws.Range("B92").Value = 4.608299999999998
$ws.Range("A93").Value = -21.27899999999999
$ws.Range("A100").Value = -21.78179999999999
